$wb = $excel.ActiveWorkbook

# ----- Summary sheet -----
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B2").Value = 0.6123595505617978
$wsSummary.Range("C2").Value = 0.5634249471458774
$wsSummary.Range("D2").Value = 0.99812734082397
$wsSummary.Range("E2").Value = 0.7202702702702702
$wsSummary.Range("F2").Value = 0.8646982478909799
$wsSummary.Range("G2").Value = 0.9693620593172916
$wsSummary.Range("H2").Value = 0.7892697330583962
$wsSummary.Range("I2").Value = 533
$wsSummary.Range("J2").Value = 413
$wsSummary.Range("K2").Value = 121
$wsSummary.Range("L2").Value = 1

# ----- Classification Report sheet -----
$wsClass = $wb.Worksheets.Item("Classification Report")
$wsClass.Range("B2").Value = 0.9918032786885246
$wsClass.Range("C2").Value = 0.2265917602996255
$wsClass.Range("D2").Value = 0.3689024390243902

$wsClass.Range("B3").Value = 0.5634249471458774
$wsClass.Range("C3").Value = 0.99812734082397
$wsClass.Range("D3").Value = 0.7202702702702702

$wsClass.Range("B4").Value = 0.6123595505617978
$wsClass.Range("C4").Value = 0.6123595505617978
$wsClass.Range("D4").Value = 0.6123595505617978
$wsClass.Range("E4").Value = 0.6123595505617978

$wsClass.Range("B5").Value = 0.777614112917201
$wsClass.Range("C5").Value = 0.6123595505617978
$wsClass.Range("D5").Value = 0.5445863546473302

$wsClass.Range("B6").Value = 0.7776141129172011
$wsClass.Range("C6").Value = 0.6123595505617978
$wsClass.Range("D6").Value = 0.5445863546473303

# ----- Confusion Matrix sheet -----
$wsConf = $wb.Worksheets.Item("Confusion Matrix")
$wsConf.Range("B2").Value = 121
$wsConf.Range("C2").Value = 413
$wsConf.Range("B3").Value = 1
$wsConf.Range("C3").Value = 533
